$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new data rows (row 7 and row 8) that were previously blank.
# The order of entry matters for how new shared strings get appended,
# so values are written in the same order they were originally typed.
$ws.Range("B7").Value = "Ferreteria 3"
$ws.Range("A7").Value = "666"
$ws.Range("C7").Value = "Melqui"
$ws.Range("D7").Value = "3176794454"
$ws.Range("E7").Value = "121212"
$ws.Range("F7").Value = "2222"

$ws.Range("B8").Value = "imagine2"
$ws.Range("A8").Value = "4545"
$ws.Range("C8").Value = "xD"
$ws.Range("D8").Value = "121212"
$ws.Range("E8").Value = "212121"
$ws.Range("F8").Value = "2121"

# Match the final active selection recorded in the workbook
$ws.Range("B7").Select()
